$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '69.865.26'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  +1.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.515.38'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -0.18%  '

$ws.Range("E4").Value = '  +0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '606.58'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +4.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '171.27'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -2.19%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.617'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -1.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.511.49'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -0.08%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.199'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +4.70%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.68'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -1.53%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.583'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  -3.17%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.22'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -0.33%  '

$ws.Range("E14").Value = '  +0.53%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.081.62'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -0.14%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '619.71'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -8.38%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.38'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.86%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.520.28'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +0.07%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.928.67'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.22%  '

$ws.Range("E20").Value = '  -2.00%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.30'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.885'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -2.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.93'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -11.98%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '15.80'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -2.88%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '95.98'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -2.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.85'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -0.44%  '

$ws.Range("E27").Value = '  +0.03%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.60'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -3.04%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.25'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.63%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '33.14'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.29%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.45'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -3.89%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.08'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -4.48%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.33'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -1.63%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.98'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -6.01%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '569.19'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.78'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -1.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.50'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.19%  '

$ws.Range("B38").Value = 'OKB'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '56.97'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -0.67%  '

$ws.Range("B39").Value = 'Hedera'
$ws.Range("C39").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.101'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.06%  '

$ws.Range("E40").Value = '  +0.10%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0452'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  +2.33%  '

$ws.Range("E42").Value = '  +2.56%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.327'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -3.84%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.332.29'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -3.13%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '33.17'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -1.36%  '

$ws.Range("E46").Value = '  +1.64%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0₃0702'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -1.09%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.62'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.53%  '

$ws.Range("E49").Value = '  -3.61%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '136.02'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  +3.07%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.67'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  +3.36%  '
